$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3030.875
$ws.Range("I125").Value = 2650
$ws.Range("J125").Value = 3665.6667
$ws.Range("K125").Value = 23850
$ws.Range("L125").Value = 32991.0003
$ws.Range("M125").Value = -21390
$ws.Range("N125").Value = -37911.0003

$ws.Range("H131").Value = 3670.6316
$ws.Range("I131").Value = 1794.8182
$ws.Range("K131").Value = 5384.4546
$ws.Range("M131").Value = -344.4546

$ws.Range("H132").Value = 1602.3214
$ws.Range("I132").Value = 1119.375
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 3358.125
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -828.125
$ws.Range("N132").Value = -18560

$ws.Range("H135").Value = 1517.9375
$ws.Range("I135").Value = 1253.4546
$ws.Range("J135").Value = 2099.8
$ws.Range("K135").Value = 11281.0914
$ws.Range("L135").Value = 18898.2
$ws.Range("M135").Value = -8746.091400000001
$ws.Range("N135").Value = -23968.2

$ws.Range("H137").Value = 5383.4863
$ws.Range("I137").Value = 5245.923
$ws.Range("J137").Value = 5708.636
$ws.Range("K137").Value = 15737.769
$ws.Range("L137").Value = 17125.908
$ws.Range("M137").Value = -13187.769
$ws.Range("N137").Value = -22225.908

$ws.Range("H138").Value = 2489.7437
$ws.Range("J138").Value = 2813.2063
$ws.Range("L138").Value = 8439.618899999999
$ws.Range("N138").Value = -18719.6189

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 456.2857
$ws.Range("I4").Value = 480.25
$ws.Range("J4").Value = 424.33334
$ws.Range("K4").Value = 480.25
$ws.Range("L4").Value = 424.33334
$ws.Range("M4").Value = -364.25
$ws.Range("N4").Value = -656.33334

$ws.Range("H28").Value = 17500
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 17500
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 17500
$ws.Range("M28").Value = ""
$ws.Range("N28").Value = -17884

$ws.Range("H31").Value = 31058.75
$ws.Range("I31").Value = 4411.6665
$ws.Range("K31").Value = 4411.6665
$ws.Range("M31").Value = -4117.6665

$ws.Range("H32").Value = 14291371
$ws.Range("I32").Value = 17858856
$ws.Range("K32").Value = 17858856
$ws.Range("M32").Value = -17858569

$ws.Range("H97").Value = 1136.75
$ws.Range("I97").Value = 930.82355
$ws.Range("K97").Value = 930.82355
$ws.Range("M97").Value = -434.82355

$ws.Range("H99").Value = 17500
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 17500
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 17500
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = -23490

$ws.Range("H108").Value = 79363.336
$ws.Range("J108").Value = 87995
$ws.Range("L108").Value = 87995
$ws.Range("N108").Value = -95675

$ws.Range("H112").Value = 28749.75
$ws.Range("J112").Value = 28749.75
$ws.Range("L112").Value = 28749.75
$ws.Range("N112").Value = -31703.75

$ws.Range("H122").Value = 4945.4443
$ws.Range("I122").Value = 4512
$ws.Range("J122").Value = 4999.625
$ws.Range("K122").Value = 13536
$ws.Range("L122").Value = 14998.875
$ws.Range("M122").Value = -11086
$ws.Range("N122").Value = -19898.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 59992.5
$ws.Range("J92").Value = 59992.5
$ws.Range("L92").Value = 59992.5
$ws.Range("N92").Value = -64984.5

$ws.Range("H100").Value = 27663
$ws.Range("J100").Value = 27663
$ws.Range("L100").Value = 27663
$ws.Range("N100").Value = -29827

$ws.Range("H134").Value = 5268144.5
$ws.Range("I134").Value = 3947.889
$ws.Range("J134").Value = 10005921
$ws.Range("K134").Value = 11843.667
$ws.Range("L134").Value = 30017763
$ws.Range("M134").Value = -9308.667000000001
$ws.Range("N134").Value = -30022833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 22919
$ws.Range("J88").Value = 22919
$ws.Range("L88").Value = 22919
$ws.Range("N88").Value = -23731

$ws.Range("H91").Value = 22919
$ws.Range("J91").Value = 22919
$ws.Range("L91").Value = 22919
$ws.Range("N91").Value = -25727

$ws.Range("H99").Value = 3779.8333
$ws.Range("I99").Value = 3568.6667
$ws.Range("J99").Value = 3991
$ws.Range("K99").Value = 3568.6667
$ws.Range("L99").Value = 3991
$ws.Range("M99").Value = -2070.6667
$ws.Range("N99").Value = -6987

$ws.Range("H122").Value = 2109.6667
$ws.Range("I122").Value = 1682.8334
$ws.Range("K122").Value = 5048.5002
$ws.Range("M122").Value = -2598.5002

$ws.Range("H126").Value = 3779.8333
$ws.Range("I126").Value = 3568.6667
$ws.Range("J126").Value = 3991
$ws.Range("K126").Value = 10706.0001
$ws.Range("L126").Value = 11973
$ws.Range("M126").Value = -8236.000100000001
$ws.Range("N126").Value = -16913

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 389.46155
$ws.Range("I34").Value = 655
$ws.Range("J34").Value = 341.18182
$ws.Range("K34").Value = 1965
$ws.Range("L34").Value = 1023.54546
$ws.Range("M34").Value = -1881
$ws.Range("N34").Value = -1191.54546

$ws.Range("H107").Value = 1562.2858
$ws.Range("J107").Value = 1739.3334
$ws.Range("L107").Value = 5218.0002
$ws.Range("N107").Value = -9058.0002

$ws.Range("H131").Value = 5462.7
$ws.Range("J131").Value = 5462.7
$ws.Range("L131").Value = 16388.1
$ws.Range("N131").Value = -26468.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 12733.333
$ws.Range("J105").Value = 12733.333
$ws.Range("L105").Value = 12733.333
$ws.Range("N105").Value = -19721.333

$ws.Range("H122").Value = 2380.3
$ws.Range("I122").Value = 1999.5
$ws.Range("J122").Value = 2951.5
$ws.Range("K122").Value = 5998.5
$ws.Range("L122").Value = 8854.5
$ws.Range("M122").Value = -3548.5
$ws.Range("N122").Value = -13754.5

$ws.Range("H126").Value = 3752.5862
$ws.Range("I126").Value = 3291.9333
$ws.Range("K126").Value = 9875.7999
$ws.Range("M126").Value = -7405.7999

$ws.Range("H132").Value = 29418976
$ws.Range("I132").Value = 38465324
$ws.Range("K132").Value = 115395972
$ws.Range("M132").Value = -115393442

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1367
$ws.Range("J22").Value = 1211.7778
$ws.Range("L22").Value = 1211.7778
$ws.Range("N22").Value = -1801.7778

$ws.Range("H27").Value = 1367
$ws.Range("J27").Value = 1211.7778
$ws.Range("L27").Value = 1211.7778
$ws.Range("N27").Value = -1425.7778

$ws.Range("H46").Value = 3984.1667
$ws.Range("J46").Value = 4975.875
$ws.Range("L46").Value = 4975.875
$ws.Range("N46").Value = -5351.875

$ws.Range("H61").Value = 2923.4
$ws.Range("I61").Value = 2654.25
$ws.Range("K61").Value = 2654.25
$ws.Range("M61").Value = -2452.25

$ws.Range("H113").Value = 2923.4
$ws.Range("I113").Value = 2654.25
$ws.Range("K113").Value = 2654.25
$ws.Range("M113").Value = -484.25

$ws.Range("H132").Value = 2023303.6
$ws.Range("I132").Value = 52254
$ws.Range("K132").Value = 156762
$ws.Range("M132").Value = -154232

$ws.Range("H136").Value = 123231.92
$ws.Range("I136").Value = 23402.2
$ws.Range("K136").Value = 70206.60000000001
$ws.Range("M136").Value = -67656.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 796
$ws.Range("I100").Value = 820.5
$ws.Range("J100").Value = 747
$ws.Range("K100").Value = 1641
$ws.Range("L100").Value = 1494
$ws.Range("M100").Value = -1100
$ws.Range("N100").Value = -2576

$ws.Range("H107").Value = 19231486
$ws.Range("I107").Value = 22727930
$ws.Range("K107").Value = 68183790
$ws.Range("M107").Value = -68181870

$ws.Range("H109").Value = 106000
$ws.Range("J109").Value = 106000
$ws.Range("L109").Value = 106000
$ws.Range("N109").Value = -108774

$ws.Range("H122").Value = 2343.0571
$ws.Range("I122").Value = 2396.7932
$ws.Range("J122").Value = 2083.3333
$ws.Range("K122").Value = 7190.3796
$ws.Range("L122").Value = 6249.999899999999
$ws.Range("M122").Value = -4740.3796
$ws.Range("N122").Value = -11149.9999

$ws.Range("H136").Value = 7792.6665
$ws.Range("J136").Value = 5414.2856
$ws.Range("L136").Value = 16242.8568
$ws.Range("N136").Value = -21342.8568
